# Update the "requirements" tracking sheet:
#  - Photo/album printing requirements (R-PD6 "print individual photo",
#    R-PD7 "print entire album") and the orphan-photo cleanup script
#    requirement (R-C) are now implemented, so mark them Complete with a
#    completion date.
#  - Sheet1 had an accidental duplicate pair of rows (R-PD8 / R-PD9);
#    remove the extra copy now that the sheet is being tidied up.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Remove the duplicated R-PD8 / R-PD9 rows (rows 21-22 duplicate the
#     later rows 30-31 exactly) ---------------------------------------------
$ws.Rows.Item(21).Resize(2).Delete()

# After the delete, the rows shift up by two:
#   R-PD6 is now row 26, R-PD7 is row 27, R-C is row 39.

# --- R-PD6 ("users can print an individual photo") -> Complete -------------
$ws.Range("B26").Value = "Complete"
$ws.Range("B26").Font.Color = 5287936   # matches the existing "Complete" green
$ws.Range("C2").Copy()                  # grab the existing date-formatted style
$ws.Range("C26").PasteSpecial(-4122)    # xlPasteFormats
$ws.Range("C26").Value2 = 40848          # 11/1/2011

# --- R-PD7 ("users can print an entire album") already had a date; just
#     confirm/leave it untouched (still Complete) -----------------------
#     (no change needed - row kept its existing 10/28/2011 completion date)

# --- R-C (anyone can run the orphaned-photo cleanup script) -> Complete ----
$ws.Range("B39").Value = "Complete"
$ws.Range("B39").Font.Color = 5287936
$ws.Range("C39").PasteSpecial(-4122)    # xlPasteFormats (reuse date style)
$ws.Range("C39").Value2 = 40848          # 11/1/2011

# --- Update the saved selection/scroll position used when the file was
#     last reviewed --------------------------------------------------------
$ws.Range("B26").Select()
